$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 54 data for 2026-01-17 run.
# Column A holds a literal date-text (matches existing rows' storage as text,
# not an Excel date serial), so force Text format before assignment, then
# drop the formatting override so no style index is left on the cell.
$dateCell = $ws.Range("A54")
$dateCell.NumberFormat = "@"
$dateCell.Value = "01/17/2026"
$dateCell.ClearFormats()

$ws.Range("B54").Value = 12724.3
$ws.Range("C54").Value = 0.22765093445746
$ws.Range("D54").Value = 0.77234906554254
$ws.Range("E54").Value = -141.28
$ws.Range("F54").Value = -20.62
$ws.Range("G54").Value = -20914.16
$ws.Range("H54").Value = -68.03
$ws.Range("I54").Value = -256.15
$ws.Range("J54").Value = -8.119999999999999
$ws.Range("K54").Value = -21170.31
$ws.Range("L54").Value = -62.46
